$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.028.20"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.207.94"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'604.67"
$ws.Range("E5").Value = "  +4.44%  "
$ws.Range("D6").Value = "'153.17"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D8").Value = "3.208.60"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "'0.507"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "'38.39"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "3.736.54"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.141.74"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'7.44"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("D18").Value = "3.206.00"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'510.66"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "'15.48"
$ws.Range("E21").Value = "  +4.08%  "
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'15.23"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'8.00"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").Value = "'85.09"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").Value = "'9.13"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").Value = "'2.85"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.78"
$ws.Range("E31").Value = "  +6.92%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'28.08"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "'1.21"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "'6.59"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").Value = "'55.26"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "'0.0902"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'478.47"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").Value = "'0.0419"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").Value = "'8.82"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("D42").Value = "'0.296"
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "2.936.80"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("D46").Value = "0.0₃0636"
$ws.Range("E46").Value = "  +3.94%  "
$ws.Range("D47").Value = "'28.70"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").Value = "'34.03"
$ws.Range("E51").Value = "  +4.47%  "
